$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.215.38'
$ws.Range("E2").Value = '  +2.46%  '

$ws.Range("D3").Value = '1.588.49'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("E4").Value = '  +1.18%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '213.11'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.93%  '

$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("E7").Value = '  +1.01%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '23.89'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +5.99%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.250'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("E10").Value = '  +0.42%  '

$ws.Range("E11").Value = '  +2.35%  '

$ws.Range("D12").Value = '1.815.95'
$ws.Range("E12").Value = '  +1.09%  '

$ws.Range("D13").Value = '1.595.46'
$ws.Range("E13").Value = '  +0.45%  '

$ws.Range("E14").Value = '  +1.25%  '

$ws.Range("E15").Value = '  -0.65%  '

$ws.Range("D16").Value = '28.262.52'
$ws.Range("E16").Value = '  +2.72%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '63.17'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.86%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '226.95'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.45%  '

$ws.Range("D19").Value = '0.0₃0707'
$ws.Range("E19").Value = '  +0.28%  '

$ws.Range("E20").Value = '  -1.17%  '

$ws.Range("E21").Value = '  +1.02%  '

$ws.Range("E22").Value = '  -1.61%  '

$ws.Range("E23").Value = '  -1.14%  '

$ws.Range("E24").Value = '  -0.27%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '151.78'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.78%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '15.17'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.25%  '

$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("E28").Value = '  -1.24%  '

$ws.Range("E29").Value = '  +1.00%  '

$ws.Range("E30").Value = '  -0.47%  '

$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("E32").Value = '  -0.60%  '

$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("D34").Value = '1.398.29'
$ws.Range("E34").Value = '  -4.18%  '

$ws.Range("E35").Value = '  -1.83%  '

$ws.Range("E36").Value = '  -8.44%  '

$ws.Range("E37").Value = '  +1.29%  '

$ws.Range("E38").Value = '  +0.13%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.54'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +6.02%  '

$ws.Range("E40").Value = '  -0.14%  '

$ws.Range("E41").Value = '  -0.57%  '

$ws.Range("E42").Value = '  +0.96%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.87'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.95%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.59'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -3.13%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.981'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.55%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '64.22'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.11%  '

$ws.Range("D47").Value = '1.725.29'
$ws.Range("E47").Value = '  +0.82%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '87.39'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.13%  '

$ws.Range("E49").Value = '  +1.85%  '

$ws.Range("E50").Value = '  +7.38%  '

$ws.Range("E51").Value = '  -0.39%  '
